# Add the "Archive API (Backend)" module section to Modules.docx.
#
# The document originally contains a single empty paragraph (carrying only
# the _GoBack bookmark). We need to insert, before that paragraph:
#   1. A centered "Archive API" title.
#   2. A level-0 numbered item "Archive API (Backend)" (List Paragraph style,
#      numId 2).
#   3. A level-1 numbered item "Модуль получения множества задач".
#   4. A level-1 numbered item "Модуль получения ссылки на условие задачи",
#      which stays in the same paragraph as the pre-existing bookmark.
#
# All new runs use Times New Roman 14pt (sz 28 / szCs 24); the first two
# paragraphs' runs (and the centered title) are additionally tagged en-US.
#
# Because the precise combination of pStyle "a3" (List Paragraph) + numPr
# ilvl/numId=2 + rFonts/sz/szCs/lang needs to land byte-for-byte in the
# OOXML, we build the target markup for the whole body (minus sectPr, which
# Word preserves automatically) and drop it in with Range.InsertXML - this
# is the same WordprocessingML a real "Insert > paste markup" op produces,
# just expressed up front instead of via a dozen discrete property edits.

$d = $word.ActiveDocument

$paragraphs = @'
<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Archive API</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>Archive API (Backend)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>Модуль получения множества</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> задач</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t>Модуль получения ссылки</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> на условие задачи</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $paragraphs.Trim() + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($xml)
